# Add age calculation and fix tenure for new hires in projection utils
# This script updates the projection results for new hires (rows 2-6,
# years 1-5) with recalculated headcount, participation, contribution
# and compensation figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B=102; C=101; D=85;  E=0.8415841584158416;  F=0.8333333333333334;
           G=0.1003363526021026; H=0.08361362716841886;
           I=460561.0260389551; J=167666.5132204776; L=167666.5132204776;
           M=628227.5392594326; N=10081246.3288; O=9673505.398699997;
           P=0.01663152627681457; Q=0.01733254971284866 }
    3 = @{ B=103; C=103; D=87;  E=0.8446601941747572;  F=0.8446601941747572;
           G=0.09903479425028895; H=0.08365074854150621;
           I=486903.624841487; J=178239.3926277735; L=178239.3926277735;
           M=665143.0174692603; N=10637203.005464; O=10229829.847461;
           P=0.01675622741581762; Q=0.01742349533526324 }
    4 = @{ B=104; C=103; D=88;  E=0.8543689320388349;  F=0.8461538461538461;
           G=0.09855902998287509; H=0.08339610229320202;
           I=508874.7839486722; J=182138.599587577; L=182138.599587577;
           M=691013.3835362492; N=10795112.11552792; O=10386767.76278483;
           P=0.01687232125413361; Q=0.01753563801052419 }
    5 = @{ B=105; C=105; D=89;  E=0.8476190476190476;  F=0.8476190476190476;
           G=0.09827781108665552; H=0.08330214463535565;
           I=531695.2573335718; J=190974.872508424; L=190974.872508424;
           M=722670.1298419957; N=11349892.00119376; O=10939197.31786837;
           P=0.01682614006268409; Q=0.01745785060449367 }
    6 = @{ B=106; C=106; D=89;  E=0.839622641509434;   F=0.839622641509434;
           G=0.0979325194061124; H=0.082226360633434;
           I=547947.509474281; J=196039.8158940278; L=196039.8158940278;
           M=743987.3253683088; N=11668294.54402957; O=11253829.02020442;
           P=0.016801068498424; Q=0.01741983244476792 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
